# Auto-generated Excel COM-interop script
# Updates market-price derived columns (H-N) on several sheets to reflect
# refreshed data from the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 12
$ws.Range("H12").Value = 1862.3334
$ws.Range("I12").Value = 2643.5
$ws.Range("K12").Value = 2643.5
$ws.Range("M12").Value = -2473.5

# Row 40
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

# Row 106
$ws.Range("H106").Value = 3056.4167
$ws.Range("I106").Value = 3098.3635
$ws.Range("K106").Value = 3098.3635
$ws.Range("M106").Value = -2467.3635

# Row 111
$ws.Range("H111").Value = 1529.8334
$ws.Range("I111").Value = 300
$ws.Range("J111").Value = 1775.8
$ws.Range("K111").Value = 900
$ws.Range("L111").Value = 5327.4
$ws.Range("M111").Value = 2167
$ws.Range("N111").Value = -11461.4

# Row 113
$ws.Range("H113").Value = 18489
$ws.Range("I113").Value = 21883.6
$ws.Range("K113").Value = 21883.6
$ws.Range("M113").Value = -18629.6

# Row 137
$ws.Range("H137").Value = 1335.875
$ws.Range("I137").Value = 1223.5555
$ws.Range("K137").Value = 3670.6665
$ws.Range("M137").Value = -1120.6665


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 2049.7827
$ws.Range("I2").Value = 1178.2858
$ws.Range("J2").Value = 3405.4443
$ws.Range("K2").Value = 1178.2858
$ws.Range("L2").Value = 3405.4443
$ws.Range("M2").Value = -1065.2858
$ws.Range("N2").Value = -3631.4443

# Row 63
$ws.Range("H63").Value = 3483.3333
$ws.Range("I63").Value = 2725
$ws.Range("K63").Value = 2725
$ws.Range("M63").Value = -2039

# Row 66
$ws.Range("H66").Value = 3483.3333
$ws.Range("I66").Value = 2725
$ws.Range("K66").Value = 13625
$ws.Range("M66").Value = -10193

# Row 110
$ws.Range("H110").Value = 1904
$ws.Range("I110").Value = 1905
$ws.Range("K110").Value = 1905
$ws.Range("M110").Value = 140

# Row 116
$ws.Range("H116").Value = 2049.7827
$ws.Range("I116").Value = 1178.2858
$ws.Range("J116").Value = 3405.4443
$ws.Range("K116").Value = 1178.2858
$ws.Range("L116").Value = 3405.4443
$ws.Range("M116").Value = 1115.7142
$ws.Range("N116").Value = -7993.4443

# Row 122
$ws.Range("H122").Value = 4125.6665
$ws.Range("I122").Value = 1800
$ws.Range("K122").Value = 5400
$ws.Range("M122").Value = -2950


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 2049.7827
$ws.Range("I3").Value = 1178.2858
$ws.Range("J3").Value = 3405.4443
$ws.Range("K3").Value = 1178.2858
$ws.Range("L3").Value = 3405.4443
$ws.Range("M3").Value = -1064.2858
$ws.Range("N3").Value = -3633.4443

# Row 86
$ws.Range("H86").Value = 2951.818
$ws.Range("I86").Value = 2108.6667
$ws.Range("J86").Value = 3963.6
$ws.Range("K86").Value = 2108.6667
$ws.Range("L86").Value = 3963.6
$ws.Range("M86").Value = -985.6667000000002
$ws.Range("N86").Value = -6209.6

# Row 89
$ws.Range("H89").Value = 2951.818
$ws.Range("I89").Value = 2108.6667
$ws.Range("J89").Value = 3963.6
$ws.Range("K89").Value = 10543.3335
$ws.Range("L89").Value = 19818
$ws.Range("M89").Value = -4927.333500000001
$ws.Range("N89").Value = -31050


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 427.86667
$ws.Range("I22").Value = 424.36365
$ws.Range("J22").Value = 437.5
$ws.Range("K22").Value = 424.36365
$ws.Range("L22").Value = 437.5
$ws.Range("M22").Value = -74.36365000000001
$ws.Range("N22").Value = -1137.5

# Row 31
$ws.Range("H31").Value = 2627.75
$ws.Range("I31").Value = 1837
$ws.Range("K31").Value = 1837
$ws.Range("M31").Value = -1542

# Row 34
$ws.Range("H34").Value = 2627.75
$ws.Range("I34").Value = 1837
$ws.Range("K34").Value = 1837
$ws.Range("M34").Value = -1635

# Row 122
$ws.Range("H122").Value = 3123.2856
$ws.Range("I122").Value = 2978
$ws.Range("J122").Value = 3995
$ws.Range("K122").Value = 8934
$ws.Range("L122").Value = 11985
$ws.Range("M122").Value = -6484
$ws.Range("N122").Value = -16885

# Row 132
$ws.Range("H132").Value = 3113.2144
$ws.Range("I132").Value = 3170.0386
$ws.Range("J132").Value = 2374.5
$ws.Range("K132").Value = 9510.1158
$ws.Range("L132").Value = 7123.5
$ws.Range("M132").Value = -6980.1158
$ws.Range("N132").Value = -12183.5


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 1117
$ws.Range("I5").Value = 1234
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 3702
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -3590
$ws.Range("N5").Value = -3224

# Row 113
$ws.Range("H113").Value = 800.3
$ws.Range("I113").Value = 736.6667
$ws.Range("J113").Value = 827.5714
$ws.Range("K113").Value = 2210.0001
$ws.Range("L113").Value = 2482.7142
$ws.Range("M113").Value = -40.0001000000002
$ws.Range("N113").Value = -6822.7142

# Row 135
$ws.Range("H135").Value = 1117
$ws.Range("I135").Value = 1234
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 11106
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -8571
$ws.Range("N135").Value = -14070


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 1082.5
$ws.Range("I102").Value = 900.875
$ws.Range("J102").Value = 1445.75
$ws.Range("K102").Value = 900.875
$ws.Range("L102").Value = 1445.75
$ws.Range("M102").Value = 721.125
$ws.Range("N102").Value = -4689.75

# Row 113
$ws.Range("H113").Value = 1377
$ws.Range("I113").Value = 1428.6666
$ws.Range("K113").Value = 1428.6666
$ws.Range("M113").Value = 741.3334

# Row 126
$ws.Range("H126").Value = 4119.4

# Row 132
$ws.Range("H132").Value = 2540.9
$ws.Range("I132").Value = 2540.9
$ws.Range("K132").Value = 7622.700000000001
$ws.Range("M132").Value = -5092.700000000001


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 46
$ws.Range("H46").Value = 4548.387
$ws.Range("I46").Value = 4466.6665
$ws.Range("J46").Value = 4625
$ws.Range("K46").Value = 4466.6665
$ws.Range("L46").Value = 4625
$ws.Range("M46").Value = -4278.6665
$ws.Range("N46").Value = -5001

# Row 61
$ws.Range("H61").Value = 1898
$ws.Range("I61").Value = 1898
$ws.Range("K61").Value = 1898
$ws.Range("M61").Value = -1696

# Row 68
$ws.Range("H68").Value = 1428.8
$ws.Range("I68").Value = 1036
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1036
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -287
$ws.Range("N68").Value = -4498

# Row 71
$ws.Range("H71").Value = 1428.8
$ws.Range("I71").Value = 1036
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 5180
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -1436
$ws.Range("N71").Value = -22488

# Row 82
$ws.Range("H82").Value = 1040
$ws.Range("I82").Value = 1200
$ws.Range("K82").Value = 1200
$ws.Range("M82").Value = -839

# Row 85
$ws.Range("H85").Value = 1040
$ws.Range("I85").Value = 1200
$ws.Range("K85").Value = 1200
$ws.Range("M85").Value = 48

# Row 113
$ws.Range("H113").Value = 1898
$ws.Range("I113").Value = 1898
$ws.Range("K113").Value = 1898
$ws.Range("M113").Value = 272

# Row 122
$ws.Range("H122").Value = 2449.5
$ws.Range("I122").Value = 2449.5
$ws.Range("K122").Value = 7348.5
$ws.Range("M122").Value = -4898.5

# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 107
$ws.Range("H107").Value = 535.63635
$ws.Range("I107").Value = 510.33334
$ws.Range("K107").Value = 1531.00002
$ws.Range("M107").Value = 388.9999800000001

# Row 113
$ws.Range("H113").Value = 440
$ws.Range("I113").Value = 440
$ws.Range("K113").Value = 1320
$ws.Range("M113").Value = 850

# Row 122
$ws.Range("H122").Value = 3750
$ws.Range("I122").Value = 3750
$ws.Range("K122").Value = 11250
$ws.Range("M122").Value = -8800

# Row 132
$ws.Range("H132").Value = 1444.7059
$ws.Range("I132").Value = 1444.7059
$ws.Range("K132").Value = 4334.1177
$ws.Range("M132").Value = -1804.1177

# Row 136
$ws.Range("H136").Value = 798.0833
$ws.Range("I136").Value = 767.56525
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 2302.69575
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = 247.3042500000001
$ws.Range("N136").Value = -9600

